# Move the "nobles" sheet to the front of the tab order, then renumber
# every sheet's name to its new zero-based position ("0", "1", "2", "3").
$wb = $excel.ActiveWorkbook

$nobles = $wb.Worksheets.Item("nobles")
$nobles.Move($wb.Worksheets.Item(1))

# Re-establish the active sheet/tab now that "nobles" leads the tab strip.
$wb.Worksheets.Item(1).Select()

# Rename every sheet to its (new) positional index, 0-based.
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).Name = [string]($i - 1)
}
